# Apply the GitHub-Actions crypto-price refresh described in the commit
# "Updated cryptos list on Sat May 20 01:01:13 UTC 2023 with GitHub Actions".
#
# Every touched cell (B/C/D/E) holds plain text in the source workbook (the
# Price/Volume columns are pre-formatted strings like "1.002" or "  +0.22%  ",
# not numbers). Excel's COM Range.Value setter auto-coerces a numeric-looking
# string (e.g. "309.25") into a real number, which would silently change the
# cell's stored type/format and corrupt values such as "6.500" -> 6.5. To keep
# these cells as text we briefly force a text NumberFormat before writing the
# value, then restore the "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '26.911.85'
Set-TextValue "E2" '  +0.06%  '

# Row 3
Set-TextValue "D3" '1.812.97'
Set-TextValue "E3" '  +0.45%  '

# Row 4
Set-TextValue "D4" '1.003'
Set-TextValue "E4" '  +0.22%  '

# Row 5
Set-TextValue "D5" '309.25'
Set-TextValue "E5" '  +0.02%  '

# Row 6
Set-TextValue "E6" '  +0.20%  '

# Row 7
Set-TextValue "D7" '0.4652'
Set-TextValue "E7" '  +0.26%  '

# Row 8
Set-TextValue "D8" '0.3660'
Set-TextValue "E8" '  -1.80%  '

# Row 9
Set-TextValue "D9" '0.07337'
Set-TextValue "E9" '  -0.28%  '

# Row 10
Set-TextValue "D10" '0.8682'
Set-TextValue "E10" '  -0.23%  '

# Row 11
Set-TextValue "D11" '20.27'
Set-TextValue "E11" '  -0.57%  '

# Row 12
Set-TextValue "D12" '1.813.08'
Set-TextValue "E12" '  -0.69%  '

# Row 13
Set-TextValue "D13" '5.354'
Set-TextValue "E13" '  -0.13%  '

# Row 14
Set-TextValue "D14" '0.07097'
Set-TextValue "E14" '  +0.94%  '

# Row 15
Set-TextValue "D15" '6.500'
Set-TextValue "E15" '  -0.25%  '

# Row 16
Set-TextValue "D16" '91.17'
Set-TextValue "E16" '  +0.36%  '

# Row 17
Set-TextValue "D17" '1.004'
Set-TextValue "E17" '  +0.27%  '

# Row 18
Set-TextValue "D18" '0.000008686'
Set-TextValue "E18" '  -0.56%  '

# Row 19
Set-TextValue "D19" '1.002'
Set-TextValue "E19" '  +0.12%  '

# Row 20
Set-TextValue "D20" '14.62'
Set-TextValue "E20" '  -0.65%  '

# Row 21
Set-TextValue "D21" '26.932.73'
Set-TextValue "E21" '  +0.06%  '

# Row 22
Set-TextValue "D22" '5.298'
Set-TextValue "E22" '  +0.24%  '

# Row 23
Set-TextValue "D23" '10.59'
Set-TextValue "E23" '  -1.10%  '

# Row 24
Set-TextValue "D24" '2.034.99'
Set-TextValue "E24" '  -0.91%  '

# Row 25
Set-TextValue "E25" '  -0.81%  '

# Row 26
Set-TextValue "D26" '150.33'
Set-TextValue "E26" '  -0.56%  '

# Row 27
Set-TextValue "D27" '2.159'
Set-TextValue "E27" '  +1.48%  '

# Row 28
Set-TextValue "D28" '18.25'
Set-TextValue "E28" '  -0.63%  '

# Row 29
Set-TextValue "D29" '5.266'
Set-TextValue "E29" '  +0.18%  '

# Row 30
Set-TextValue "D30" '115.45'
Set-TextValue "E30" '  -0.21%  '

# Row 31
Set-TextValue "D31" '0.08936'
Set-TextValue "E31" '  +0.39%  '

# Row 32
Set-TextValue "D32" '0.7556'
Set-TextValue "E32" '  -0.88%  '

# Row 33
Set-TextValue "D33" '1.156'
Set-TextValue "E33" '  +0.10%  '

# Row 34
Set-TextValue "D34" '4.474'
Set-TextValue "E34" '  +0.42%  '

# Row 35
Set-TextValue "D35" '2.919'
Set-TextValue "E35" '  +0.28%  '

# Row 36
Set-TextValue "D36" '1.003'
Set-TextValue "E36" '  +0.27%  '

# Row 37
Set-TextValue "D37" '1.084'
Set-TextValue "E37" '  -1.97%  '

# Row 38
Set-TextValue "B38" 'Hedera'
Set-TextValue "C38" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D38" '0.05274'
Set-TextValue "E38" '  +0.90%  '

# Row 39
Set-TextValue "B39" 'VeChain'
Set-TextValue "C39" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D39" '0.01952'
Set-TextValue "E39" '  -0.10%  '

# Row 40
Set-TextValue "D40" '2.972'
Set-TextValue "E40" '  +2.48%  '

# Row 41
Set-TextValue "B41" 'TheSandbox'
Set-TextValue "C41" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D41" '0.5301'
Set-TextValue "E41" '  -0.29%  '

# Row 42
Set-TextValue "B42" 'FraxShare'
Set-TextValue "C42" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D42" '7.176'
Set-TextValue "E42" '  -1.37%  '

# Row 43
Set-TextValue "D43" '2.296'
Set-TextValue "E43" '  -3.46%  '

# Row 44
Set-TextValue "D44" '0.1646'
Set-TextValue "E44" '  -0.61%  '

# Row 45
Set-TextValue "D45" '8.370'
Set-TextValue "E45" '  -1.76%  '

# Row 46
Set-TextValue "D46" '0.4852'
Set-TextValue "E46" '  -3.49%  '

# Row 47
Set-TextValue "D47" '10.35'
Set-TextValue "E47" '  +0.84%  '

# Row 48
Set-TextValue "E48" '  +0.25%  '

# Row 49
Set-TextValue "D49" '1.658'
Set-TextValue "E49" '  -0.21%  '

# Row 50
Set-TextValue "D50" '102.62'
Set-TextValue "E50" '  -0.69%  '

# Row 51
Set-TextValue "D51" '0.06290'
Set-TextValue "E51" '  -0.07%  '

